$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.027123
$ws.Range("H2").Value2 = 0.081369
$ws.Range("I2").Value2 = 0.0960827240265261
$ws.Range("J2").Value2 = 0.09608272402652611
$ws.Range("M2").Value2 = 0.3655573333333333
$ws.Range("N2").Value2 = 1.096672
$ws.Range("O2").Value2 = 0.2375364113931583
$ws.Range("P2").Value2 = 0.2572707166041556
$ws.Range("Q2").Value2 = 0.009915011551999997
$ws.Range("R2").Value2 = 0.08923510396799998
$ws.Range("S2").Value2 = 0.0228231454621402
$ws.Range("T2").Value2 = 0.02471927126358369
$ws.Range("G3").Value2 = 0.027123
$ws.Range("H3").Value2 = 0.081369
$ws.Range("I3").Value2 = 0.0960827240265261
$ws.Range("J3").Value2 = 0.09608272402652611
$ws.Range("M3").Value2 = 0.6295006666666667
$ws.Range("O3").Value2 = 0.4090448082825151
$ws.Range("P3").Value2 = 0.4430278723705731
$ws.Range("Q3").Value2 = 0.017073946582
$ws.Range("R3").Value2 = 0.153665519238
$ws.Range("S3").Value2 = 0.03930213942869217
$ws.Range("T3").Value2 = 0.04256732479704081
$ws.Range("G4").Value2 = 0.027123
$ws.Range("H4").Value2 = 0.081369
$ws.Range("I4").Value2 = 0.0960827240265261
$ws.Range("J4").Value2 = 0.09608272402652611
$ws.Range("M4").Value2 = 0.1127876666666667
$ws.Range("N4").Value2 = 0.338363
$ws.Range("O4").Value2 = 0.07328857923629238
$ws.Range("P4").Value2 = 0.07937732656831935
$ws.Range("Q4").Value2 = 0.003059139883
$ws.Range("R4").Value2 = 0.027532258947
$ws.Range("S4").Value2 = 0.007041766333056872
$ws.Range("T4").Value2 = 0.007626789762627267
$ws.Range("G5").Value2 = 0.027123
$ws.Range("H5").Value2 = 0.081369
$ws.Range("I5").Value2 = 0.0960827240265261
$ws.Range("J5").Value2 = 0.09608272402652611
$ws.Range("M5").Value2 = 0.3541425
$ws.Range("N5").Value2 = 0.7082850000000001
$ws.Range("O5").Value2 = 0.2301191383708208
$ws.Range("P5").Value2 = 0.1661581489360305
$ws.Range("Q5").Value2 = 0.0096054070275
$ws.Range("R5").Value2 = 0.057632442165
$ws.Range("S5").Value2 = 0.02211047366530555
$ws.Range("T5").Value2 = 0.01596492756897904
$ws.Range("G6").Value2 = 0.027123
$ws.Range("H6").Value2 = 0.081369
$ws.Range("I6").Value2 = 0.0960827240265261
$ws.Range("J6").Value2 = 0.09608272402652611
$ws.Range("M6").Value2 = 0.07696466666666667
$ws.Range("N6").Value2 = 0.230894
$ws.Range("O6").Value2 = 0.05001106271721345
$ws.Range("P6").Value2 = 0.0541659355209214
$ws.Range("Q6").Value2 = 0.002087512654
$ws.Range("R6").Value2 = 0.018787613886
$ws.Range("S6").Value2 = 0.004805199137331308
$ws.Range("T6").Value2 = 0.005204410634295299
$ws.Range("I7").Value2 = 0.9039172759734738
$ws.Range("J7").Value2 = 0.9039172759734738
$ws.Range("M7").Value2 = 0.3655573333333333
$ws.Range("N7").Value2 = 1.096672
$ws.Range("O7").Value2 = 0.2375364113931583
$ws.Range("P7").Value2 = 0.2572707166041556
$ws.Range("Q7").Value2 = 0.09327743695999997
$ws.Range("R7").Value2 = 0.8394969326399998
$ws.Range("S7").Value2 = 0.2147132659310181
$ws.Range("T7").Value2 = 0.2325514453405719
$ws.Range("I8").Value2 = 0.9039172759734738
$ws.Range("J8").Value2 = 0.9039172759734738
$ws.Range("M8").Value2 = 0.6295006666666667
$ws.Range("O8").Value2 = 0.4090448082825151
$ws.Range("P8").Value2 = 0.4430278723705731
$ws.Range("S8").Value2 = 0.3697426688538229
$ws.Range("T8").Value2 = 0.4004605475735323
$ws.Range("I9").Value2 = 0.9039172759734738
$ws.Range("J9").Value2 = 0.9039172759734738
$ws.Range("M9").Value2 = 0.1127876666666667
$ws.Range("N9").Value2 = 0.338363
$ws.Range("O9").Value2 = 0.07328857923629238
$ws.Range("P9").Value2 = 0.07937732656831935
$ws.Range("Q9").Value2 = 0.02877946496499999
$ws.Range("R9").Value2 = 0.259015184685
$ws.Range("S9").Value2 = 0.0662468129032355
$ws.Range("T9").Value2 = 0.07175053680569207
$ws.Range("I10").Value2 = 0.9039172759734738
$ws.Range("J10").Value2 = 0.9039172759734738
$ws.Range("M10").Value2 = 0.3541425
$ws.Range("N10").Value2 = 0.7082850000000001
$ws.Range("O10").Value2 = 0.2301191383708208
$ws.Range("P10").Value2 = 0.1661581489360305
$ws.Range("Q10").Value2 = 0.0903647710125
$ws.Range("R10").Value2 = 0.542188626075
$ws.Range("S10").Value2 = 0.2080086647055152
$ws.Range("T10").Value2 = 0.1501932213670514
$ws.Range("I11").Value2 = 0.9039172759734738
$ws.Range("J11").Value2 = 0.9039172759734738
$ws.Range("M11").Value2 = 0.07696466666666667
$ws.Range("N11").Value2 = 0.230894
$ws.Range("O11").Value2 = 0.05001106271721345
$ws.Range("P11").Value2 = 0.0541659355209214
$ws.Range("Q11").Value2 = 0.01963868917
$ws.Range("R11").Value2 = 0.17674820253
$ws.Range("S11").Value2 = 0.04520586357988214
$ws.Range("T11").Value2 = 0.0489615248866261